# Update import price list: remove the "local_name" column (column C)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire column C ("local_name"), shifting subsequent columns left
$ws.Range("C:C").Delete()

# Reset the selection to A2, matching the post-edit cursor position
$ws.Range("A2").Select()
